$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.520183
$ws.Range("H2").Value = 1.560549
$ws.Range("I2").Value = 0.03656880080220595
$ws.Range("J2").Value = 0.03656880080220595
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 122.328922
$ws.Range("N2").Value = 366.986766
$ws.Range("O2").Value = 0.9783373008518612
$ws.Range("P2").Value = 0.9783373008518613
$ws.Range("Q2").Value = 63.63342563272599
$ws.Range("R2").Value = 572.7008306945339
$ws.Range("S2").Value = 0.03577662187221955
$ws.Range("T2").Value = 0.03577662187221956

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.520183
$ws.Range("H3").Value = 1.560549
$ws.Range("I3").Value = 0.03656880080220595
$ws.Range("J3").Value = 0.03656880080220595
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3863573333333334
$ws.Range("N3").Value = 1.159072
$ws.Range("O3").Value = 0.003089929874945324
$ws.Range("P3").Value = 0.003089929874945324
$ws.Range("Q3").Value = 0.2009765167253333
$ws.Range("R3").Value = 1.808788650528
$ws.Range("S3").Value = 0.0001129950300896607
$ws.Range("T3").Value = 0.0001129950300896607

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.520183
$ws.Range("H4").Value = 1.560549
$ws.Range("I4").Value = 0.03656880080220595
$ws.Range("J4").Value = 0.03656880080220595
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.322294
$ws.Range("N4").Value = 6.966882000000001
$ws.Range("O4").Value = 0.0185727692731934
$ws.Range("P4").Value = 0.0185727692731934
$ws.Range("Q4").Value = 1.208017859802
$ws.Range("R4").Value = 10.872160738218
$ws.Range("S4").Value = 0.0006791838998967411
$ws.Range("T4").Value = 0.0006791838998967411

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.984906
$ws.Range("H5").Value = 38.954718
$ws.Range("I5").Value = 0.912837291778795
$ws.Range("J5").Value = 0.9128372917787949
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 122.328922
$ws.Range("N5").Value = 366.986766
$ws.Range("O5").Value = 0.9783373008518612
$ws.Range("P5").Value = 0.9783373008518613
$ws.Range("Q5").Value = 1588.429553251332
$ws.Range("R5").Value = 14295.86597926199
$ws.Range("S5").Value = 0.8930627721557892
$ws.Range("T5").Value = 0.8930627721557892

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 12.984906
$ws.Range("H6").Value = 38.954718
$ws.Range("I6").Value = 0.912837291778795
$ws.Range("J6").Value = 0.9128372917787949
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3863573333333334
$ws.Range("N6").Value = 1.159072
$ws.Range("O6").Value = 0.003089929874945324
$ws.Range("P6").Value = 0.003089929874945324
$ws.Range("Q6").Value = 5.016813655744001
$ws.Range("R6").Value = 45.15132290169601
$ws.Range("S6").Value = 0.00282060321883148
$ws.Range("T6").Value = 0.002820603218831479

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 12.984906
$ws.Range("H7").Value = 38.954718
$ws.Range("I7").Value = 0.912837291778795
$ws.Range("J7").Value = 0.9128372917787949
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.322294
$ws.Range("N7").Value = 6.966882000000001
$ws.Range("O7").Value = 0.0185727692731934
$ws.Range("P7").Value = 0.0185727692731934
$ws.Range("Q7").Value = 30.15476929436401
$ws.Range("R7").Value = 271.392923649276
$ws.Range("S7").Value = 0.01695391640417429
$ws.Range("T7").Value = 0.01695391640417428

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.719687
$ws.Range("H8").Value = 2.159061
$ws.Range("I8").Value = 0.05059390741899907
$ws.Range("J8").Value = 0.05059390741899907
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 122.328922
$ws.Range("N8").Value = 366.986766
$ws.Range("O8").Value = 0.9783373008518612
$ws.Range("P8").Value = 0.9783373008518613
$ws.Range("Q8").Value = 88.03853488741399
$ws.Range("R8").Value = 792.3468139867259
$ws.Range("S8").Value = 0.04949790682385251
$ws.Range("T8").Value = 0.04949790682385252

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.719687
$ws.Range("H9").Value = 2.159061
$ws.Range("I9").Value = 0.05059390741899907
$ws.Range("J9").Value = 0.05059390741899907
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.3863573333333334
$ws.Range("N9").Value = 1.159072
$ws.Range("O9").Value = 0.003089929874945324
$ws.Range("P9").Value = 0.003089929874945324
$ws.Range("Q9").Value = 0.2780563501546667
$ws.Range("R9").Value = 2.502507151392
$ws.Range("S9").Value = 0.0001563316260241831
$ws.Range("T9").Value = 0.0001563316260241831

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.719687
$ws.Range("H10").Value = 2.159061
$ws.Range("I10").Value = 0.05059390741899907
$ws.Range("J10").Value = 0.05059390741899907
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.322294
$ws.Range("N10").Value = 6.966882000000001
$ws.Range("O10").Value = 0.0185727692731934
$ws.Range("P10").Value = 0.0185727692731934
$ws.Range("Q10").Value = 1.671324801978
$ws.Range("R10").Value = 15.041923217802
$ws.Range("S10").Value = 0.0009396689691223778
$ws.Range("T10").Value = 0.0009396689691223778
